{"js": "// Update the worksheet date and every division problem's two operands.\n// The document is a \"date header\" paragraph followed by a table of\n// division problems (\"NN\u00f7N=\" strings). Each old string is unique in the\n// document, so a scoped search()+insertText(Replace) round-trip is\n// enough to retarget a specific cell without disturbing its run\n// formatting (font/size/paragraph alignment).\n//\n// NOTE on ordering: a couple of the new values collide with *other*\n// cells' old values (e.g. one cell's old text becomes a text that used\n// to be a different cell's old text). If we replaced in naive order, a\n// later search for that original text would match BOTH the untouched\n// cell and the text we just wrote, corrupting the wrong cell. To avoid\n// that, replacements are applied in an order where a rule only runs\n// after any other rule whose \"replace\" text equals this rule's \"find\"\n// text has already completed.\nconst replacements = [\n  { find: \"2023-09-04 Monday\", replace: \"2023-09-05 Tuesday\" },\n  { find: \"58\u00f72=\", replace: \"12\u00f75=\" },\n  { find: \"72\u00f77=\", replace: \"82\u00f74=\" },\n  { find: \"62\u00f78=\", replace: \"98\u00f78=\" },\n  { find: \"36\u00f74=\", replace: \"25\u00f77=\" },\n  { find: \"36\u00f73=\", replace: \"78\u00f76=\" },\n  { find: \"76\u00f79=\", replace: \"57\u00f76=\" },\n  { find: \"41\u00f75=\", replace: \"28\u00f79=\" },\n  { find: \"17\u00f79=\", replace: \"34\u00f77=\" },\n  { find: \"81\u00f76=\", replace: \"53\u00f77=\" },\n  { find: \"34\u00f74=\", replace: \"27\u00f76=\" },\n  { find: \"82\u00f78=\", replace: \"93\u00f77=\" },\n  { find: \"13\u00f74=\", replace: \"55\u00f78=\" },\n  { find: \"89\u00f74=\", replace: \"76\u00f73=\" },\n  { find: \"49\u00f72=\", replace: \"66\u00f76=\" },\n  { find: \"93\u00f72=\", replace: \"32\u00f73=\" },\n  { find: \"35\u00f77=\", replace: \"28\u00f73=\" },\n  { find: \"25\u00f76=\", replace: \"45\u00f76=\" },\n  { find: \"17\u00f74=\", replace: \"54\u00f74=\" },\n  { find: \"93\u00f79=\", replace: \"54\u00f77=\" },\n  { find: \"96\u00f76=\", replace: \"45\u00f79=\" },\n  { find: \"31\u00f74=\", replace: \"84\u00f72=\" },\n  { find: \"47\u00f75=\", replace: \"29\u00f74=\" },\n  { find: \"42\u00f76=\", replace: \"75\u00f75=\" },\n  { find: \"19\u00f73=\", replace: \"51\u00f72=\" },\n  { find: \"95\u00f74=\", replace: \"42\u00f76=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every division problem's two operands.\n# The document is a \"date header\" paragraph followed by a table of\n# division problems (\"NN\u00f7N=\" strings). Each old string is unique in the\n# document, so Find/Replace scoped to a fresh $d.Content range each time\n# retargets exactly one run without disturbing its formatting\n# (font/size/paragraph alignment survive because only the <w:t> text is\n# touched).\n#\n# NOTE on ordering: a couple of the new values collide with *other*\n# cells' old values (e.g. one cell's old text becomes a text that used\n# to be a different cell's old text). If we replaced in naive order, a\n# later Find for that original text would match BOTH the untouched\n# cell and the text we just wrote, and ReplaceAll would corrupt the\n# wrong cell. To avoid that, replacements are applied in an order where\n# a pair only runs after any other pair whose Replace text equals this\n# pair's Find text has already completed.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Find = \"2023-09-04 Monday\"; Replace = \"2023-09-05 Tuesday\" },\n    @{ Find = \"58\u00f72=\"; Replace = \"12\u00f75=\" },\n    @{ Find = \"72\u00f77=\"; Replace = \"82\u00f74=\" },\n    @{ Find = \"62\u00f78=\"; Replace = \"98\u00f78=\" },\n    @{ Find = \"36\u00f74=\"; Replace = \"25\u00f77=\" },\n    @{ Find = \"36\u00f73=\"; Replace = \"78\u00f76=\" },\n    @{ Find = \"76\u00f79=\"; Replace = \"57\u00f76=\" },\n    @{ Find = \"41\u00f75=\"; Replace = \"28\u00f79=\" },\n    @{ Find = \"17\u00f79=\"; Replace = \"34\u00f77=\" },\n    @{ Find = \"81\u00f76=\"; Replace = \"53\u00f77=\" },\n    @{ Find = \"34\u00f74=\"; Replace = \"27\u00f76=\" },\n    @{ Find = \"82\u00f78=\"; Replace = \"93\u00f77=\" },\n    @{ Find = \"13\u00f74=\"; Replace = \"55\u00f78=\" },\n    @{ Find = \"89\u00f74=\"; Replace = \"76\u00f73=\" },\n    @{ Find = \"49\u00f72=\"; Replace = \"66\u00f76=\" },\n    @{ Find = \"93\u00f72=\"; Replace = \"32\u00f73=\" },\n    @{ Find = \"35\u00f77=\"; Replace = \"28\u00f73=\" },\n    @{ Find = \"25\u00f76=\"; Replace = \"45\u00f76=\" },\n    @{ Find = \"17\u00f74=\"; Replace = \"54\u00f74=\" },\n    @{ Find = \"93\u00f79=\"; Replace = \"54\u00f77=\" },\n    @{ Find = \"96\u00f76=\"; Replace = \"45\u00f79=\" },\n    @{ Find = \"31\u00f74=\"; Replace = \"84\u00f72=\" },\n    @{ Find = \"47\u00f75=\"; Replace = \"29\u00f74=\" },\n    @{ Find = \"42\u00f76=\"; Replace = \"75\u00f75=\" },\n    @{ Find = \"19\u00f73=\"; Replace = \"51\u00f72=\" },\n    @{ Find = \"95\u00f74=\"; Replace = \"42\u00f76=\" }\n)\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $null = $range.Find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n"}
